# Regenerate merged AHB file: rows 26-59 on the single worksheet get the
# same "merged/regenerated" formatting that rows 2-25 already have:
#   - On each segment-group header row, every cell A:V is restyled to the
#     shared "header" look (grey fill + border + wrap, bold in column B).
#   - On every row in the block (header rows included), the "Änderung"
#     column (L) loses its "ÄNDERUNG" label/style and becomes a plain,
#     empty, centered grey cell - matching column L on rows 2-25.
#
# Rows 2-25 already contain the exact target formatting, so we copy
# formats from there (PasteSpecial formats-only) instead of hand-rolling
# style indices - this reuses the existing cellXfs entries instead of
# minting new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Segment-group header rows (first row of each "#" group) in the 26-59 block.
$headerRows = @(26, 30, 37, 41, 48, 55, 57)

# Every other (non-header) row in the 26-59 block.
$nonHeaderRows = @(27, 28, 29, 31, 32, 33, 34, 35, 36, 38, 39, 40, 42, 43, 44, 45, 46, 47, 49, 50, 51, 52, 53, 54, 56, 58, 59)

# Template rows already in the target format.
$headerTemplate = $ws.Range("A2:V2")
$changeColTemplate = $ws.Range("L2")

foreach ($r in $headerRows) {
    $headerTemplate.Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)
}

foreach ($r in $nonHeaderRows) {
    $changeColTemplate.Copy()
    $ws.Range("L" + $r).PasteSpecial(-4122)
}

# Drop the old "ÄNDERUNG" label text from column L for the whole block
# (header rows' L cell is already blank after the paste above, but
# ClearContents is harmless/idempotent there).
foreach ($r in ($headerRows + $nonHeaderRows)) {
    $ws.Range("L" + $r).ClearContents()
}

$excel.CutCopyMode = $false
